# Auto-generated edit script to apply precision-rounding updates
# as described in the commit diff for test_results.xlsx

$wb = $excel.ActiveWorkbook

$wsInput = $wb.Worksheets.Item("Input Parameters")
$wsOutput = $wb.Worksheets.Item("Output Results")

# --- Updates on "Input Parameters" sheet (Latitude/Longitude/Yield rounding) ---
$wsInput.Range("B2").Value = 24.1
$wsInput.Range("C2").Value = 90.41
$wsInput.Range("B3").Value = 24.1
$wsInput.Range("C3").Value = 90.41
$wsInput.Range("B4").Value = 24.1
$wsInput.Range("C4").Value = 90.41
$wsInput.Range("B5").Value = 24.1
$wsInput.Range("C5").Value = 90.41
$wsInput.Range("B6").Value = 24.1
$wsInput.Range("C6").Value = 90.41
$wsInput.Range("B7").Value = 24.1
$wsInput.Range("C7").Value = 90.41
$wsInput.Range("B8").Value = 24.1
$wsInput.Range("C8").Value = 90.41
$wsInput.Range("B9").Value = 24.92
$wsInput.Range("C9").Value = 89.94
$wsInput.Range("B10").Value = 24.92
$wsInput.Range("C10").Value = 89.94
$wsInput.Range("B11").Value = 24.92
$wsInput.Range("C11").Value = 89.94
$wsInput.Range("B12").Value = 24.92
$wsInput.Range("C12").Value = 89.94
$wsInput.Range("B13").Value = 24.92
$wsInput.Range("C13").Value = 89.94
$wsInput.Range("M13").Value = 35.18
$wsInput.Range("B14").Value = 24.92
$wsInput.Range("C14").Value = 89.94
$wsInput.Range("M14").Value = 32.66
$wsInput.Range("B15").Value = 24.1
$wsInput.Range("C15").Value = 90.41
$wsInput.Range("B16").Value = 23.49
$wsInput.Range("C16").Value = 89.42
$wsInput.Range("B17").Value = 22.7
$wsInput.Range("C17").Value = 90.37

# --- Updates on "Output Results" sheet (Yield / Seasonal irrigation recalculated values) ---
$wsOutput.Range("F2").Value = 8.252035301011807
$wsOutput.Range("F3").Value = 8.265331847337418
$wsOutput.Range("F4").Value = 8.331373152267158
$wsOutput.Range("F5").Value = 8.329031300172565
$wsOutput.Range("F6").Value = 8.252035301011807
$wsOutput.Range("F7").Value = 8.265331847337418
$wsOutput.Range("F8").Value = 8.331373152267158
$wsOutput.Range("F9").Value = 8.329031300172565
$wsOutput.Range("F10").Value = 8.252035301011807
$wsOutput.Range("F11").Value = 8.265331847337418
$wsOutput.Range("F12").Value = 8.331373152267158
$wsOutput.Range("F13").Value = 8.329031300172565
$wsOutput.Range("F14").Value = 5.799594592429121
$wsOutput.Range("F15").Value = 5.724272415628232
$wsOutput.Range("F16").Value = 5.425883008718891
$wsOutput.Range("G16").Value = 575
$wsOutput.Range("F17").Value = 5.799594592429121
$wsOutput.Range("F18").Value = 5.724272415628232
$wsOutput.Range("F19").Value = 5.425883008718891
$wsOutput.Range("G19").Value = 575
$wsOutput.Range("F20").Value = 5.799594592429121
$wsOutput.Range("F21").Value = 5.724272415628232
$wsOutput.Range("F22").Value = 5.425883008718891
$wsOutput.Range("G22").Value = 575
$wsOutput.Range("F23").Value = 5.799594592429121
$wsOutput.Range("F24").Value = 5.724272415628232
$wsOutput.Range("F25").Value = 5.425883008718891
$wsOutput.Range("G25").Value = 575
$wsOutput.Range("F26").Value = 2.199050389369717
$wsOutput.Range("G26").Value = 324.6294128218001
$wsOutput.Range("F27").Value = 3.619606235631393
$wsOutput.Range("F28").Value = 1.494692662904003
$wsOutput.Range("F29").Value = 2.199050389369717
$wsOutput.Range("G29").Value = 324.6294128218001
$wsOutput.Range("F30").Value = 3.619606235631393
$wsOutput.Range("F31").Value = 1.494692662904003
$wsOutput.Range("F32").Value = 1.75629891650748
$wsOutput.Range("F33").Value = 1.404878119678743
$wsOutput.Range("F34").Value = 1.775297960730243
$wsOutput.Range("F35").Value = 1.75629891650748
$wsOutput.Range("F36").Value = 1.404878119678743
$wsOutput.Range("F37").Value = 1.775297960730243
$wsOutput.Range("F38").Value = 6.219102643765745
$wsOutput.Range("F39").Value = 2.667879848391237
$wsOutput.Range("F40").Value = 1.910364104782561
$wsOutput.Range("F41").Value = 6.219102643765745
$wsOutput.Range("F42").Value = 2.667879848391237
$wsOutput.Range("F43").Value = 1.910364104782561
$wsOutput.Range("G44").Value = 1166.947191939173
$wsOutput.Range("G45").Value = 1200.390170187957
$wsOutput.Range("G46").Value = 1150.102147786943
$wsOutput.Range("G47").Value = 722.8977391671099
$wsOutput.Range("F48").Value = 7.015346729218368
$wsOutput.Range("G48").Value = 689.147681439479
$wsOutput.Range("F49").Value = 7.043660632884428
$wsOutput.Range("G49").Value = 687.9765759548468
$wsOutput.Range("G50").Value = 674.9028563351791
$wsOutput.Range("G51").Value = 696.9296587599061
$wsOutput.Range("F52").Value = 7.015211013062076
$wsOutput.Range("G52").Value = 698.2311625269382
$wsOutput.Range("F53").Value = 7.047831954663301
$wsOutput.Range("G53").Value = 665.171884249018
$wsOutput.Range("F54").Value = 7.08299744934908
$wsOutput.Range("G54").Value = 674.0966184177332
